# "ajustado comentarios dos arquivos" - fill in the "Tipo" header label in
# column K (row 7) and the missing "Linha" numbers in column J for the
# rows that were left blank.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K7").Value = "Tipo"

$ws.Range("J9").Value = 12
$ws.Range("J10").Value = 13
$ws.Range("J11").Value = 16
$ws.Range("J12").Value = 17
$ws.Range("J15").Value = 18
$ws.Range("J16").Value = 19
$ws.Range("J33").Value = 24
$ws.Range("J36").Value = 8
$ws.Range("J37").Value = 11
